# JournalDeTravail.xlsx - "Work diary and documentation" commit
#
# Appends six new work-log entries (rows 97-102) to the "Tableau1" table on
# the only worksheet, covering documentation/review work performed on
# 2022-05-30 and 2022-05-31 (Excel serials 44711 / 44712).
#
# Growing the table via ListRows.Add() first makes Excel own the job of
# extending the table ref / AutoFilter range (and the sheet's used range);
# the cell values are then written directly so every new cell lands with
# the same "wrap text" column style the rest of the sheet already uses.
# The H96 total (=SUM(C:C)) recalculates automatically once the new
# "Duree (heures)" values are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tableau1")

$rowsToAdd = 6
for ($i = 0; $i -lt $rowsToAdd; $i++) {
    $null = $lo.ListRows.Add()
}

# Firstrow of the newly appended block (Tableau1 started at A1:F96, i.e.
# 95 data rows below the header -> first new data row is sheet row 97).
$firstNewRow = 97

$entries = @(
    @{ Date = 44711; Type = "Réalisation"; Duree = 0.5;  Description = "Revue de la documentation";                              Remarque = "Travail a la maison" },
    @{ Date = 44711; Type = "Réalisation"; Duree = 0.5;  Description = "Documentation des erreurs restantes";                    Remarque = "Beaucoup d'erreur. Recherche sur comment les résoudres" },
    @{ Date = 44711; Type = "Réalisation"; Duree = 0.25; Description = "Spellcheck rapide de la documentation";                  Remarque = $null },
    @{ Date = 44711; Type = "Réalisation"; Duree = 1;    Description = "Documentation de la conclusion";                         Remarque = "Suite possible au projet, réparation de la mise en page" },
    @{ Date = 44712; Type = "Réalisation"; Duree = 1.25; Description = "Documentation - Guide de mise en service";               Remarque = $null },
    @{ Date = 44712; Type = "Réalisation"; Duree = 0.25; Description = "Revue de la documentation une dernière fois ce soir";     Remarque = $null }
)

for ($i = 0; $i -lt $entries.Count; $i++) {
    $r = $firstNewRow + $i
    $entry = $entries[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $entry.Date
    $dateCell.NumberFormat = "dd/mm/yyyy"

    $ws.Cells.Item($r, 2).Value = $entry.Type
    $ws.Cells.Item($r, 3).Value = $entry.Duree
    $ws.Cells.Item($r, 4).Value = $entry.Description

    if ($entry.Remarque) {
        $ws.Cells.Item($r, 5).Value = $entry.Remarque
    }
}

# Mirror the saved selection/scroll state left behind by the author.
$null = $ws.Range("E102").Select()

Write-Host "Added $rowsToAdd work-diary rows ($firstNewRow..$($firstNewRow + $entries.Count - 1))."
